$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.410.04'
$ws.Range('E2').Value = '  -1.51%  '

$ws.Range('D3').Value = '1.643.28'
$ws.Range('E3').Value = '  -0.65%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.14%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.002'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '298.78'
$ws.Range('E6').Value = '  -1.79%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3784'
$ws.Range('E7').Value = '  -0.87%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3528'
$ws.Range('E8').Value = '  -2.35%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '50.06'
$ws.Range('E9').Value = '  -2.10%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08069'
$ws.Range('E10').Value = '  -1.62%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.211'

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.003'
$ws.Range('E12').Value = '  +0.12%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.03'
$ws.Range('E13').Value = '  -2.90%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.378'
$ws.Range('E14').Value = '  -2.40%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.306'
$ws.Range('E15').Value = '  -1.68%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001198'
$ws.Range('E16').Value = '  -3.24%  '

$ws.Range('D17').Value = '1.643.53'
$ws.Range('E17').Value = '  -0.38%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '96.87'
$ws.Range('E18').Value = '  -1.02%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06971'
$ws.Range('E19').Value = '  -0.03%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.744'
$ws.Range('E20').Value = '  -0.51%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.37'
$ws.Range('E21').Value = '  -2.13%  '

$ws.Range('E22').Value = '  +0.01%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.38'
$ws.Range('E23').Value = '  -2.60%  '

$ws.Range('D24').Value = '23.423.63'
$ws.Range('E24').Value = '  -1.42%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.511'
$ws.Range('E25').Value = '  -1.58%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.890'
$ws.Range('E26').Value = '  -6.18%  '

$ws.Range('E27').Value = '  -2.17%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '152.53'
$ws.Range('E28').Value = '  +1.04%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.197'
$ws.Range('E29').Value = '  -0.31%  '

$ws.Range('E30').Value = '  -1.78%  '

$ws.Range('D31').Value = '1.824.64'
$ws.Range('E31').Value = '  +0.04%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.900'
$ws.Range('E32').Value = '  -0.05%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.142'
$ws.Range('E33').Value = '  +0.23%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.40'
$ws.Range('E34').Value = '  -4.52%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9853'

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02689'
$ws.Range('E36').Value = '  -5.21%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.08730'
$ws.Range('E37').Value = '  -1.20%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2435'
$ws.Range('E38').Value = '  -3.35%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.909'
$ws.Range('E39').Value = '  -3.88%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06787'
$ws.Range('E40').Value = '  -5.35%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '12.91'
$ws.Range('E41').Value = '  -1.40%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6855'
$ws.Range('E42').Value = '  -3.06%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.287'
$ws.Range('E43').Value = '  -4.20%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.59'
$ws.Range('E44').Value = '  -2.32%  '

$ws.Range('E45').Value = '  +0.13%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6331'
$ws.Range('E46').Value = '  -3.25%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.248'
$ws.Range('E47').Value = '  -3.51%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.901'
$ws.Range('E48').Value = '  -1.58%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07722'
$ws.Range('E49').Value = '  -3.22%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '127.34'
$ws.Range('E50').Value = '  -1.01%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.145'
$ws.Range('E51').Value = '  -4.13%  '
